$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97

$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "Balanço Geral"
$ws.Cells.Item($row, 3).Value = "Iluminação"
$ws.Cells.Item($row, 4).Value = "2025-04-09T13:21"
$ws.Cells.Item($row, 5).Value = "Negativo"
$ws.Cells.Item($row, 6).Value = "Iluminação na Ponte Leonel Brizola incomoda moradores próximos ao loca. Repórter *ao vivo*. Ponte ganhou iluminação rosa. Alegam que tem inclusive atrapalhado o sono deles. Equipe procurou a prefeitura para saber se existe possibilidades de representar. *com nota da Prefeitura*"
